# Swap the contents of columns G and H (header + all data rows).
# Before: G = "Best Explanation (flipped)" label column, H = "P(H|O)" probability column
# After:  G = "P(H|O)" probability column, H = "Best Explanation (flipped)" label column
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$gRange = $ws.Range("G1:G$lastRow")
$hRange = $ws.Range("H1:H$lastRow")

$gValues = $gRange.Value()
$hValues = $hRange.Value()

$gRange.Value = $hValues
$hRange.Value = $gValues
